$d = $word.ActiveDocument

$d.Content.Find.Execute("14×57=798", $true, $false, $false, $false, $false, $true, 1, $false, "39×93=3627", 2) | Out-Null
$d.Content.Find.Execute("73×57=4161", $true, $false, $false, $false, $false, $true, 1, $false, "83×71=5893", 2) | Out-Null
$d.Content.Find.Execute("31×62=1922", $true, $false, $false, $false, $false, $true, 1, $false, "77×84=6468", 2) | Out-Null
$d.Content.Find.Execute("80×20=1600", $true, $false, $false, $false, $false, $true, 1, $false, "38×37=1406", 2) | Out-Null
$d.Content.Find.Execute("77×60=4620", $true, $false, $false, $false, $false, $true, 1, $false, "34×16=544", 2) | Out-Null
$d.Content.Find.Execute("85×96=8160", $true, $false, $false, $false, $false, $true, 1, $false, "98×64=6272", 2) | Out-Null
$d.Content.Find.Execute("16×24=384", $true, $false, $false, $false, $false, $true, 1, $false, "38×52=1976", 2) | Out-Null
$d.Content.Find.Execute("92×17=1564", $true, $false, $false, $false, $false, $true, 1, $false, "51×70=3570", 2) | Out-Null
$d.Content.Find.Execute("76×14=1064", $true, $false, $false, $false, $false, $true, 1, $false, "60×54=3240", 2) | Out-Null
$d.Content.Find.Execute("87×89=7743", $true, $false, $false, $false, $false, $true, 1, $false, "46×75=3450", 2) | Out-Null
$d.Content.Find.Execute("14×62=868", $true, $false, $false, $false, $false, $true, 1, $false, "88×81=7128", 2) | Out-Null
$d.Content.Find.Execute("43×99=4257", $true, $false, $false, $false, $false, $true, 1, $false, "72×79=5688", 2) | Out-Null
$d.Content.Find.Execute("65×67=4355", $true, $false, $false, $false, $false, $true, 1, $false, "85×64=5440", 2) | Out-Null
$d.Content.Find.Execute("93×88=8184", $true, $false, $false, $false, $false, $true, 1, $false, "57×46=2622", 2) | Out-Null
$d.Content.Find.Execute("65×57=3705", $true, $false, $false, $false, $false, $true, 1, $false, "55×77=4235", 2) | Out-Null
$d.Content.Find.Execute("93×30=2790", $true, $false, $false, $false, $false, $true, 1, $false, "74×54=3996", 2) | Out-Null
$d.Content.Find.Execute("40×14=560", $true, $false, $false, $false, $false, $true, 1, $false, "49×90=4410", 2) | Out-Null
$d.Content.Find.Execute("65×47=3055", $true, $false, $false, $false, $false, $true, 1, $false, "60×69=4140", 2) | Out-Null
$d.Content.Find.Execute("70×39=2730", $true, $false, $false, $false, $false, $true, 1, $false, "85×11=935", 2) | Out-Null
$d.Content.Find.Execute("90×63=5670", $true, $false, $false, $false, $false, $true, 1, $false, "59×83=4897", 2) | Out-Null
$d.Content.Find.Execute("14×14=196", $true, $false, $false, $false, $false, $true, 1, $false, "43×71=3053", 2) | Out-Null
$d.Content.Find.Execute("57×87=4959", $true, $false, $false, $false, $false, $true, 1, $false, "97×95=9215", 2) | Out-Null
$d.Content.Find.Execute("63×49=3087", $true, $false, $false, $false, $false, $true, 1, $false, "75×42=3150", 2) | Out-Null
$d.Content.Find.Execute("35×37=1295", $true, $false, $false, $false, $false, $true, 1, $false, "60×46=2760", 2) | Out-Null
$d.Content.Find.Execute("73×70=5110", $true, $false, $false, $false, $false, $true, 1, $false, "60×12=720", 2) | Out-Null

Write-Host "Replacements applied."
